$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column T: "Customer Reference" header + "Test3" value, matching the
# existing header row's style (reuse the same style as the other headers).
$ws.Range("T1").Value = "Customer Reference"
$ws.Range("T2").Value = "Test3"

$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

$ws.Range("T2").Select()
